# Restructure the interpolate test workbook:
#   - insert a new "about" sheet in front, describing the test file
#   - rename "data_need_interpolation" -> "dummy_input"
#   - append a new, empty "output" sheet at the end (becomes the active tab)

$wb = $excel.ActiveWorkbook

# 1) New "about" sheet, inserted before the current first sheet ("path").
$firstSheet = $wb.Worksheets.Item(1)
$aboutSheet = $wb.Worksheets.Add($firstSheet)
$aboutSheet.Name = "about"

$aboutSheet.Range("B2").Value = "This file is used to test"
$aboutSheet.Range("B3").Value = "interpolation.m"
$aboutSheet.Range("B4").Value = "It has a dummy data table table with timeseries missing some years"
$aboutSheet.Range("B5").Value = "the interpolate scripts adds the missing years and creats a output table with more rows"
[void]$aboutSheet.Range("B6").Select()

# 2) Rename the second data sheet and move the selection off its old cell.
$dummySheet = $wb.Worksheets.Item("data_need_interpolation")
$dummySheet.Name = "dummy_input"
[void]$dummySheet.Range("E12").Select()

# 3) New "output" sheet appended after the last sheet; this becomes active.
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$outputSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$outputSheet.Name = "output"
[void]$outputSheet.Range("H16").Select()
